$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values (columns A and B, rows 1-3) and add two new rows (4-5)
$ws.Range("A1").Value = 0.087974784335598757
$ws.Range("B1").Value = -0.087974785859716451

$ws.Range("A2").Value = -0.028603004232630083
$ws.Range("B2").Value = 0.028603002562140486

$ws.Range("A3").Value = -0.0021519143033771728
$ws.Range("B3").Value = 0.0021519127695541224

$ws.Range("A4").Value = 0.012012106307611056
$ws.Range("B4").Value = -0.012012107973542895

$ws.Range("A5").Value = -0.0086647554732060972
$ws.Range("B5").Value = 0.0086647537959538262

# Widen columns A and B by one character's width (matches diff target widths)
$ws.Columns.Item(1).ColumnWidth = 14.592447916666666
$ws.Columns.Item(2).ColumnWidth = 13.592447916666666
